# Update Handback status timestamps as part of "Generate Report for Handback".
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the e089686e row (row 4, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-17 14:46:02"

# zh-cn sheet: Correspond Handoff/Handback DateTime for the e089686e row (row 4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-17 14:45:55"
$wsZhCn.Range("K4").Value = "2016-08-17 14:46:32"

# de-de sheet: Correspond Handback DateTime for the e089686e row (row 4)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-17 14:46:40"
